# Apply scheduled-runner profit recalculations across all sheets.
# Generated from the authoritative cell-level diff: each (sheet,row) block
# updates a subset of columns H..N with the commit's recomputed values.

$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# --- ALC sheet ---
Set-RowValues "ALC" 64 @{ "H" = 71434440; "J" = 7666.3335; "L" = 7666.3335; "N" = -8162.3335 }
Set-RowValues "ALC" 67 @{ "H" = 71434440; "J" = 7666.3335; "L" = 7666.3335; "N" = -9382.333500000001 }
Set-RowValues "ALC" 70 @{ "H" = 4800.353; "I" = 5606.4287; "J" = 4236.1; "K" = 16819.2861; "L" = 12708.3; "M" = -16549.2861; "N" = -13248.3 }
Set-RowValues "ALC" 73 @{ "H" = 4800.353; "I" = 5606.4287; "J" = 4236.1; "K" = 16819.2861; "L" = 12708.3; "M" = -15883.2861; "N" = -14580.3 }
Set-RowValues "ALC" 107 @{ "H" = 1797; "I" = 2685; "J" = 1353; "K" = 2685; "L" = 1353; "M" = -765; "N" = -5193 }
Set-RowValues "ALC" 111 @{ "H" = 1822.6428; "I" = 1020.4; "K" = 3061.2; "M" = 5.800000000000182 }
Set-RowValues "ALC" 112 @{ "H" = 3885.84; "J" = 4517.3; "L" = 13551.9; "N" = -15767.9 }
Set-RowValues "ALC" 129 @{ "H" = 1547.8; "J" = 3175; "L" = 9525; "N" = -19525 }
Set-RowValues "ALC" 137 @{ "H" = 9127.322; "I" = 2990.1333; "K" = 8970.3999; "M" = -6420.3999 }
Set-RowValues "ALC" 138 @{ "H" = 1918.16; "J" = 2132.8442; "L" = 6398.5326; "N" = -16678.5326 }
# --- ARM sheet ---
Set-RowValues "ARM" 32 @{ "H" = 2984.5078; "I" = 2223.6316; "K" = 2223.6316; "M" = -1936.6316 }
Set-RowValues "ARM" 45 @{ "H" = 3846.5; "I" = 3639.875; "K" = 3639.875; "M" = -3262.875 }
Set-RowValues "ARM" 102 @{ "H" = 3196.9375; "I" = 2410.0667; "K" = 2410.0667; "M" = -788.0666999999999 }
Set-RowValues "ARM" 122 @{ "H" = 3685.3914; "I" = 1590.7333; "J" = 7612.875; "K" = 4772.199900000001; "L" = 22838.625; "M" = -2322.199900000001; "N" = -27738.625 }
# --- BSM sheet ---
Set-RowValues "BSM" 20 @{ "H" = 18041.836; "I" = 5441.68; "K" = 5441.68; "M" = -5194.68 }
Set-RowValues "BSM" 86 @{ "H" = 2446.4194; "I" = 2515.2; "J" = 2321.3635; "K" = 2515.2; "L" = 2321.3635; "M" = -1392.2; "N" = -4567.363499999999 }
Set-RowValues "BSM" 89 @{ "H" = 2446.4194; "I" = 2515.2; "J" = 2321.3635; "K" = 12576; "L" = 11606.8175; "M" = -6960; "N" = -22838.8175 }
Set-RowValues "BSM" 94 @{ "H" = 1417.6852; "J" = 2531.1765; "L" = 2531.1765; "N" = -3433.1765 }
# --- CRP sheet ---
Set-RowValues "CRP" 31 @{ "H" = 56139.383; "I" = 45396; "J" = 68673.336; "K" = 45396; "L" = 68673.336; "M" = -45101; "N" = -69263.336 }
Set-RowValues "CRP" 34 @{ "H" = 56139.383; "I" = 45396; "J" = 68673.336; "K" = 45396; "L" = 68673.336; "M" = -45194; "N" = -69077.336 }
Set-RowValues "CRP" 99 @{ "H" = 159533.83; "I" = 267465.94; "J" = 8428.9; "K" = 267465.94; "L" = 8428.9; "M" = -265967.94; "N" = -11424.9 }
Set-RowValues "CRP" 107 @{ "H" = 1064.2667; "I" = 1351.9375; "J" = 735.5; "K" = 1351.9375; "L" = 735.5; "M" = 568.0625; "N" = -4575.5 }
Set-RowValues "CRP" 122 @{ "H" = 3230.5833; "I" = 2632; "K" = 7896; "M" = -5446 }
Set-RowValues "CRP" 126 @{ "H" = 159533.83; "I" = 267465.94; "J" = 8428.9; "K" = 802397.8200000001; "L" = 25286.7; "M" = -799927.8200000001; "N" = -30226.7 }
Set-RowValues "CRP" 132 @{ "H" = 2600.5; "I" = 2617.3333; "J" = 2499.5; "K" = 7851.999899999999; "L" = 7498.5; "M" = -5321.999899999999; "N" = -12558.5 }
Set-RowValues "CRP" 134 @{ "H" = 443107.12; "I" = 1959.6552; "J" = 1722434.8; "K" = 5878.9656; "L" = 5167304.4; "M" = -3343.9656; "N" = -5172374.4 }
# --- CUL sheet ---
Set-RowValues "CUL" 26 @{ "H" = 2047.8 }
# --- GSM sheet ---
Set-RowValues "GSM" 80 @{ "H" = 9895.242; "I" = 8539.434999999999; "J" = 13013.6; "K" = 8539.434999999999; "L" = 13013.6; "M" = -7541.434999999999; "N" = -15009.6 }
Set-RowValues "GSM" 83 @{ "H" = 9895.242; "I" = 8539.434999999999; "J" = 13013.6; "K" = 42697.175; "L" = 65068; "M" = -37705.175; "N" = -75052 }
Set-RowValues "GSM" 97 @{ "H" = 1331.931; "I" = 1299.8334; "K" = 1299.8334; "M" = -803.8334 }
Set-RowValues "GSM" 122 @{ "H" = 4540.905; "I" = 5805.636; "J" = 3149.7; "K" = 17416.908; "L" = 9449.099999999999; "M" = -14966.908; "N" = -14349.1 }
Set-RowValues "GSM" 126 @{ "H" = 9165.044; "I" = 15182.75; "K" = 45548.25; "M" = -43078.25 }
# --- LTW sheet ---
Set-RowValues "LTW" 7 @{ "H" = 8950.8125; "J" = 7999.75; "L" = 7999.75; "N" = -8223.75 }
Set-RowValues "LTW" 61 @{ "H" = 3224.1943; "I" = 2508.0356; "J" = 5730.75; "K" = 2508.0356; "L" = 5730.75; "M" = -2306.0356; "N" = -6134.75 }
Set-RowValues "LTW" 113 @{ "H" = 3224.1943; "I" = 2508.0356; "J" = 5730.75; "K" = 2508.0356; "L" = 5730.75; "M" = -338.0356000000002; "N" = -10070.75 }
Set-RowValues "LTW" 122 @{ "H" = 20839600; "I" = 43484268; "J" = 6505.16; "K" = 130452804; "L" = 19515.48; "M" = -130450354; "N" = -24415.48 }
Set-RowValues "LTW" 126 @{ "H" = 8950.8125; "J" = 7999.75; "L" = 23999.25; "N" = -28939.25 }
# --- WVR sheet ---
Set-RowValues "WVR" 22 @{ "H" = 12572.571; "I" = 7004; "J" = 14800; "K" = 7004; "L" = 14800; "M" = -6711; "N" = -15386 }
Set-RowValues "WVR" 107 @{ "H" = 1248.1818; "I" = 1248.1818; "K" = 3744.5454; "M" = -1824.5454 }
Set-RowValues "WVR" 122 @{ "H" = 447970.75; "I" = 558268.7; "K" = 1674806.1; "M" = -1672356.1 }
